# chore: adapt column header formatting to respective input file names
#
#   "<name>_old" -> "<name>_FV2310"
#   "<name>_new" -> "<name>_FV2404"
#
# ... and turn the A1:U65 data range into a real Excel Table ("Table1"),
# with the header row frozen in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row (row 1) --------------------------------------
# Column headers currently end in "_old" (left/"FV2310" block) or "_new"
# (right/"FV2404" block); the "diff" column in between is left untouched.
$oldSuffix = "_old"
$newSuffix = "_FV2310"
$oldSuffix2 = "_new"
$newSuffix2 = "_FV2404"

$lastCol = 21   # column U
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $header = [string]$cell.Value2
    if ($header.EndsWith($oldSuffix)) {
        $cell.Value = $header.Substring(0, $header.Length - $oldSuffix.Length) + $newSuffix
    } elseif ($header.EndsWith($oldSuffix2)) {
        $cell.Value = $header.Substring(0, $header.Length - $oldSuffix2.Length) + $newSuffix2
    }
}

# --- 2. Freeze the header row -----------------------------------------------
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the used range into an Excel Table ------------------------------
$dataRange = $ws.Range("A1:U65")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"
